# input: Week 4 matches
# Fill in the Week 4 match results on the "2024-Match" sheet (rows 14-17)
# and switch the active sheet / selection to reflect where the user was
# working, matching the sheet that was previously active (2023-Match).

$wb = $excel.ActiveWorkbook

$wsMatch2024 = $wb.Worksheets.Item("2024-Match")
$wsMatch2023 = $wb.Worksheets.Item("2023-Match")

# --- Row 14 ---
$wsMatch2024.Range("A14").Value = "Week 4"
$wsMatch2024.Range("B14").Value = "July"
$wsMatch2024.Range("C14").Value = 9
$wsMatch2024.Range("D14").Value = "Tuesday"
$wsMatch2024.Range("E14").Value = 0.80902777777777779
$wsMatch2024.Range("F14").Value = "APEX Charters Lone Pine Brewing"
$wsMatch2024.Range("G14").Value = "Baxter Pines FC"
$wsMatch2024.Range("H14").Value = 0
$wsMatch2024.Range("I14").Value = 1
$wsMatch2024.Range("J14").Value = "Dean Zoulamis"
$wsMatch2024.Range("K14").Value = "Patrik Udeh"

# --- Row 15 ---
$wsMatch2024.Range("A15").Value = "Week 4"
$wsMatch2024.Range("B15").Value = "July"
$wsMatch2024.Range("C15").Value = 9
$wsMatch2024.Range("D15").Value = "Tuesday"
$wsMatch2024.Range("E15").Value = 0.88194444444444453
$wsMatch2024.Range("F15").Value = "Carlos Auto Repair"
$wsMatch2024.Range("G15").Value = "The Escape Room"
$wsMatch2024.Range("H15").Value = 1
$wsMatch2024.Range("I15").Value = 1
$wsMatch2024.Range("J15").Value = "Patrik Udeh"
$wsMatch2024.Range("K15").Value = "Dean Zoulamis"

# --- Row 16 ---
$wsMatch2024.Range("A16").Value = "Week 4"
$wsMatch2024.Range("B16").Value = "July"
$wsMatch2024.Range("C16").Value = 11
$wsMatch2024.Range("D16").Value = "Thursday"
$wsMatch2024.Range("E16").Value = 0.80902777777777779
$wsMatch2024.Range("F16").Value = "Old Port FC"
$wsMatch2024.Range("G16").Value = "Farmers FC"
$wsMatch2024.Range("H16").Value = 0
$wsMatch2024.Range("I16").Value = 1
$wsMatch2024.Range("J16").Value = "Caleb Lamb"
$wsMatch2024.Range("K16").Value = "Eric"

# --- Row 17 ---
$wsMatch2024.Range("A17").Value = "Week 4"
$wsMatch2024.Range("B17").Value = "July"
$wsMatch2024.Range("C17").Value = 11
$wsMatch2024.Range("D17").Value = "Thursday"
$wsMatch2024.Range("E17").Value = 0.88194444444444453
$wsMatch2024.Range("F17").Value = "One Love FC"
$wsMatch2024.Range("G17").Value = "Thunder FC"
$wsMatch2024.Range("H17").Value = 4
$wsMatch2024.Range("I17").Value = 5
$wsMatch2024.Range("J17").Value = "Eric"
$wsMatch2024.Range("K17").Value = "Caleb Lamb"

# --- View / selection state ---
# "2023-Match" keeps its own selection (unchanged from before the edit).
[void]$wsMatch2023.Range("J26").Select()

# Previously "2023-Match" was the active tab; now "2024-Match" becomes the
# active tab, with the selection parked at J22.
[void]$wsMatch2024.Activate()
[void]$wsMatch2024.Range("J22").Select()
